$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = '@'
    $rng.Value = $value
    $rng.Style = $origStyle
}

Set-TextValue 'D2' '29.899.16'
Set-TextValue 'E2' '  +0.17%  '
Set-TextValue 'D3' '1.889.78'
Set-TextValue 'E3' '  -0.08%  '
Set-TextValue 'E4' '  +0.06%  '
Set-TextValue 'D5' '0.7691'
Set-TextValue 'E5' '  -1.52%  '
Set-TextValue 'E6' '  -0.56%  '
Set-TextValue 'E7' '  +0.03%  '
Set-TextValue 'D8' '0.3136'
Set-TextValue 'E8' '  -0.24%  '
Set-TextValue 'D9' '25.69'
Set-TextValue 'E9' '  +1.20%  '
Set-TextValue 'D10' '0.07138'
Set-TextValue 'E10' '  -2.39%  '
Set-TextValue 'E11' '  +5.06%  '
Set-TextValue 'D12' '0.7646'
Set-TextValue 'E12' '  -0.26%  '
Set-TextValue 'D13' '1.917.22'
Set-TextValue 'E13' '  +2.15%  '
Set-TextValue 'D14' '5.373'
Set-TextValue 'E14' '  -1.72%  '
Set-TextValue 'D15' '93.75'
Set-TextValue 'E15' '  +0.65%  '
Set-TextValue 'D16' '6.145'
Set-TextValue 'E16' '  -0.98%  '
Set-TextValue 'D17' '29.924.89'
Set-TextValue 'E17' '  +0.31%  '
Set-TextValue 'E18' '  -1.07%  '
Set-TextValue 'D19' '244.48'
Set-TextValue 'E19' '  -0.52%  '
Set-TextValue 'D20' '0.000007818'
Set-TextValue 'E20' '  -0.72%  '
Set-TextValue 'E21' '  -0.03%  '
Set-TextValue 'D22' '8.028'
Set-TextValue 'E22' '  -1.51%  '
Set-TextValue 'D23' '1.000'
Set-TextValue 'E23' '  +0.08%  '
Set-TextValue 'D24' '0.1630'
Set-TextValue 'E24' '  +2.75%  '
Set-TextValue 'E25' '  -0.71%  '
Set-TextValue 'D26' '163.08'
Set-TextValue 'E26' '  +0.81%  '
Set-TextValue 'D28' '2.041'
Set-TextValue 'E28' '  +0.26%  '
Set-TextValue 'D29' '1.517'
Set-TextValue 'E29' '  +4.47%  '
Set-TextValue 'D30' '1.536'
Set-TextValue 'E30' '  -0.39%  '
Set-TextValue 'D31' '4.511'
Set-TextValue 'E31' '  +0.80%  '
Set-TextValue 'D33' '0.05455'
Set-TextValue 'D34' '1.243'
Set-TextValue 'D35' '0.7465'
Set-TextValue 'E35' '  -1.18%  '
Set-TextValue 'E36' '  +0.31%  '
Set-TextValue 'D37' '2.698'
Set-TextValue 'E37' '  +2.21%  '
Set-TextValue 'D38' '0.01951'
Set-TextValue 'E38' '  +0.84%  '
Set-TextValue 'D39' '2.780'
Set-TextValue 'E39' '  -0.07%  '
Set-TextValue 'D40' '0.4476'
Set-TextValue 'E40' '  +0.56%  '
Set-TextValue 'D41' '1.102.84'
Set-TextValue 'E41' '  -3.39%  '
Set-TextValue 'D42' '73.30'
Set-TextValue 'E43' '  +1.86%  '
Set-TextValue 'D44' '0.8520'
Set-TextValue 'E44' '  -0.36%  '
Set-TextValue 'E45' '  +0.02%  '
Set-TextValue 'D46' '103.05'
Set-TextValue 'E46' '  +1.15%  '
Set-TextValue 'D47' '1.871'
Set-TextValue 'E47' '  -1.59%  '
Set-TextValue 'D48' '7.680'
Set-TextValue 'E48' '  +2.24%  '
Set-TextValue 'D49' '3.042'
Set-TextValue 'E49' '  -2.04%  '
Set-TextValue 'D50' '2.030.67'
Set-TextValue 'E50' '  +0.10%  '
Set-TextValue 'D51' '0.06083'
Set-TextValue 'E51' '  +0.30%  '
